# Reorder the categorized "RESEARCH" sheet entries alphabetically within each YEAR group.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESEARCH")

# Target values (row -> classification, count) reflecting the new alphabetical ordering.
$values = @{
    2  = @("Book Chapter", 1)
    3  = @("Conference Presentation", 2)
    4  = @("Journal Publication-Indexed", 1)
    5  = @("Mou s", 6)
    6  = @("Seminar Organized", 1)
    7  = @("Mou s", 1)
    8  = @("Book Chapter", 1)
    9  = @("Conference Attended", 2)
    10 = @("Conference Keynote", 1)
    11 = @("Conference Presentation", 4)
    12 = @("Conference Publication", 2)
    13 = @("Conference Session Chair", 1)
    14 = @("Journal Publication-Indexed", 5)
    15 = @("Journal Publication-Non Indexed", 3)
    16 = @("Mou s", 1)
    17 = @("Patent Filed", 3)
    18 = @("Seminar Organized", 1)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}
